$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "kingson@virogreenusa.com"
$ws.Range("A2").Value = "info@horizontechnology.com"
$ws.Range("A3").Value = "sales@horizontechnology.com"
$ws.Range("A4").Value = "info@ctdi.com"
